$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ52872080",
    "summ52999940",
    "summ53128037",
    "summ53351661",
    "summ53542786",
    "summ53701851",
    "summ53864273",
    "summ54022303",
    "summ54197243",
    "summ54356758",
    "summ54530042",
    "summ54691586",
    "summ54866966",
    "summ55039183",
    "summ55205730",
    "summ55367770",
    "summ55542811",
    "summ55707113",
    "summ55883178",
    "summ56059730",
    "summ56232843",
    "summ56420408",
    "summ56667161",
    "summ56844182",
    "summ57025947",
    "summ57202976",
    "summ57390216",
    "summ57565246",
    "summ57752042",
    "summ57928112",
    "summ58093648",
    "summ58260130",
    "summ58424693",
    "summ58608795",
    "summ58780965",
    "summ58954066",
    "summ59123133",
    "summ59369222",
    "summ59519693",
    "summ59689338",
    "summ59863357",
    "summ00036420",
    "summ00329226",
    "summ00487617",
    "summ00656059",
    "summ00824608",
    "summ00986849",
    "summ01139917",
    "summ01307365",
    "summ01469393"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
